# The 월요일 (Monday) sheet's 7th period row (row 8) previously stored the
# start/end times for period 7 as real numeric time-of-day values
# (07:20 / 07:25) formatted with a "h:mm" number format.
#
# The edit replaces those with literal text values "19:20" / "19:25"
# (the actual 7th-period class time), stored as plain Text-formatted
# cells - matching how every other day-of-week sheet already stores its
# period times (as shared-string text, not as numeric time serials).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("월요일")

$rng = $ws.Range("B8:C8")
$rng.NumberFormat = "@"

$ws.Range("B8").Value = "19:20"
$ws.Range("C8").Value = "19:25"
